{"js": "// Replace the multiplication problems in the document with their updated values.\n// The mapping below reflects each old expression -> new expression from the diff.\nconst replacements = [\n  [\"111\u00d75=\", \"670\u00d78=\"],\n  [\"841\u00d78=\", \"435\u00d74=\"],\n  [\"680\u00d75=\", \"651\u00d72=\"],\n  [\"274\u00d77=\", \"218\u00d77=\"],\n  [\"942\u00d76=\", \"835\u00d75=\"],\n  [\"149\u00d78=\", \"265\u00d75=\"],\n  [\"770\u00d78=\", \"222\u00d77=\"],\n  [\"619\u00d75=\", \"477\u00d78=\"],\n  [\"872\u00d79=\", \"371\u00d78=\"],\n  [\"579\u00d76=\", \"486\u00d75=\"],\n  [\"855\u00d73=\", \"688\u00d78=\"],\n  [\"936\u00d73=\", \"655\u00d76=\"],\n  [\"732\u00d78=\", \"192\u00d75=\"],\n  [\"696\u00d76=\", \"454\u00d76=\"],\n  [\"612\u00d78=\", \"542\u00d77=\"],\n  [\"969\u00d72=\", \"786\u00d76=\"],\n  [\"662\u00d73=\", \"244\u00d74=\"],\n  [\"860\u00d75=\", \"892\u00d78=\"],\n  [\"698\u00d76=\", \"410\u00d74=\"],\n  [\"672\u00d79=\", \"446\u00d72=\"],\n  [\"925\u00d75=\", \"277\u00d72=\"],\n  [\"178\u00d72=\", \"236\u00d78=\"],\n  [\"335\u00d78=\", \"307\u00d72=\"],\n  [\"721\u00d75=\", \"335\u00d73=\"],\n  [\"155\u00d74=\", \"292\u00d77=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the multiplication problems in the document with their updated values.\n# The mapping below reflects each old expression -> new expression from the diff.\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"111\u00d75=\", \"670\u00d78=\"),\n    @(\"841\u00d78=\", \"435\u00d74=\"),\n    @(\"680\u00d75=\", \"651\u00d72=\"),\n    @(\"274\u00d77=\", \"218\u00d77=\"),\n    @(\"942\u00d76=\", \"835\u00d75=\"),\n    @(\"149\u00d78=\", \"265\u00d75=\"),\n    @(\"770\u00d78=\", \"222\u00d77=\"),\n    @(\"619\u00d75=\", \"477\u00d78=\"),\n    @(\"872\u00d79=\", \"371\u00d78=\"),\n    @(\"579\u00d76=\", \"486\u00d75=\"),\n    @(\"855\u00d73=\", \"688\u00d78=\"),\n    @(\"936\u00d73=\", \"655\u00d76=\"),\n    @(\"732\u00d78=\", \"192\u00d75=\"),\n    @(\"696\u00d76=\", \"454\u00d76=\"),\n    @(\"612\u00d78=\", \"542\u00d77=\"),\n    @(\"969\u00d72=\", \"786\u00d76=\"),\n    @(\"662\u00d73=\", \"244\u00d74=\"),\n    @(\"860\u00d75=\", \"892\u00d78=\"),\n    @(\"698\u00d76=\", \"410\u00d74=\"),\n    @(\"672\u00d79=\", \"446\u00d72=\"),\n    @(\"925\u00d75=\", \"277\u00d72=\"),\n    @(\"178\u00d72=\", \"236\u00d78=\"),\n    @(\"335\u00d78=\", \"307\u00d72=\"),\n    @(\"721\u00d75=\", \"335\u00d73=\"),\n    @(\"155\u00d74=\", \"292\u00d77=\")\n)\n\nforeach ($pair in $pairs) {\n    $old = $pair[0]\n    $new = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $old\n    $find.Replacement.Text = $new\n    $find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2)\n}\n"}
